# "Add files via upload" — append two new rows (13 & 14) of product-entry
# data to Sheet1 and drop the explicit (General) number-format style that
# used to sit on A11/A12 (it now effectively "moves" onto the new A14 row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# A11 / A12 previously carried an explicit "General" number format style
# (s="2"). In the edited workbook that style attribute is gone, so reset
# both cells back to the workbook's default (Normal) style.
# ---------------------------------------------------------------------
$ws.Range("A11").Style = "Normal"
$ws.Range("A12").Style = "Normal"

# ---------------------------------------------------------------------
# Row 13 — new record
# ---------------------------------------------------------------------
$ws.Range("A13").Value = 62455
$ws.Range("B13").Value = "lym00987"
$ws.Range("C13").Value = "bhgty"
# Leading apostrophe forces this date-shaped string to stay literal text
# (matching the existing Date of Build column, which stores text, not
# real dates); re-applying the Normal style afterwards clears the
# "quote prefix" flag Excel stamps on while keeping the text untouched.
$ws.Range("D13").Value = "'2025-03-07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = 0.006
$ws.Range("F13").Value = 0.008
$ws.Range("G13").Value = 0.009
$ws.Range("H13").Value = 0.006
$ws.Range("I13").Value = 7.7
$ws.Range("J13").Value = 9.5
$ws.Range("K13").Value = 5.6
$ws.Range("L13").Value = 2.3

# ---------------------------------------------------------------------
# Row 14 — new record (A14 keeps the explicit General-format style, like
# A11/A12 used to)
# ---------------------------------------------------------------------
$ws.Range("A14").Value = 6543
$ws.Range("A14").NumberFormat = "General"
$ws.Range("B14").Value = "tz800789"
$ws.Range("C14").Value = "kmnhj"
$ws.Range("D14").Value = "'2025-03-05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = 0.004
$ws.Range("F14").Value = 0.002
$ws.Range("G14").Value = 0.008
$ws.Range("H14").Value = 0.009
$ws.Range("I14").Value = 3.4
$ws.Range("J14").Value = 5.6
$ws.Range("K14").Value = 7.8
$ws.Range("L14").Value = 9

# Selection moved down two rows to B17 (tracking the two new data rows).
$ws.Range("B17").Select()
